$wb = $excel.ActiveWorkbook

$wsTestCases = $wb.Worksheets.Item("Test Cases")
$wsCourse    = $wb.Worksheets.Item("Complete_Course_1")

# --- "Complete_Course_1" sheet updates -------------------------------------
# Widen column A to fit the new (longer) email address.
$wsCourse.Columns.Item(1).ColumnWidth = 42.307291666666664

# Username cell (A2) gets a new test user's email address. The cell keeps its
# existing (hyperlink) style; only the displayed/shared-string text changes.
$wsCourse.Range("A2").Value = "bhabani.shankar105@weboapps.com"

# Move this sheet's selection/cursor (no longer the active tab once we
# activate "Test Cases" below).
$wsCourse.Range("A13").Select() | Out-Null

# --- "Test Cases" sheet becomes the active tab/selection -------------------
$wsTestCases.Activate() | Out-Null
$wsTestCases.Range("C6").Select() | Out-Null
